# Burndown Chart.xlsx edits
# 1. Rename the sheet from "Sheet1" to "Burndown Chart"
# 2. Update the BackLog IDs legend text from "1,2,3,4,5" to "1,2,3,4"
# 3. Mark 0.2 of backlog items 1 and 4 as completed on day 5 (column G)
#    which updates the Real Burndown computed values accordingly
# 4. Move the active selection to K16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Burndown Chart"

$ws.Range("A15").Value = "1,2,3,4"
$ws.Range("A16").Value = "1,2,3,4"
$ws.Range("A17").Value = "1,2,3,4"
$ws.Range("A18").Value = "1,2,3,4"
$ws.Range("A19").Value = "1,2,3,4"

$ws.Range("G3").Value = 0.2
$ws.Range("G6").Value = 0.2

$ws.Range("K16").Select()
